# Update the Equipment sheet data per the new data set, and reset the
# active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

# New row values: Tag, Description, PumpType, PumpDriverType, DesignTemp,
# DesignPressure, Capacity, SpecificGravity, DifferentialPressure
$rows = @(
    @{ Row = 2; A = "Equip-001"; B = "DESC-9"; C = "PT-8"; D = "PDT-3"; E = 4; F = 6; G = 3; H = 3; I = 2 },
    @{ Row = 3; A = "Equip-002"; B = "DESC-3"; C = "PT-7"; D = "PDT-9"; E = 4; F = 7; G = 4; H = 7; I = 9 },
    @{ Row = 4; A = "Equip-003"; B = "DESC-2"; C = "PT-7"; D = "PDT-3"; E = 8; F = 7; G = 7; H = 2; I = 8 },
    @{ Row = 5; A = "Equip-004"; B = "DESC-3"; C = "PT-3"; D = "PDT-8"; E = 2; F = 4; G = 5; H = 3; I = 3 }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H
    $ws.Range("I" + $r.Row).Value = $r.I
}

# Move/reset the selected cell in the Equipment sheet to A2
$ws.Activate()
$ws.Range("A2").Select()
